$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells in this sheet are stored as text (inlineStr) in the
# source workbook. Force the number format to Text ("@") before writing
# so Excel does not silently reinterpret numeric-looking strings (e.g.
# "1.0000", "0.000007734", "2.655") as numbers and strip formatting.

$ws.Range('D2:E2').NumberFormat = "@"
$ws.Range('D2').Value = '30.860.34'
$ws.Range('E2').Value = '  -1.29%  '

$ws.Range('D3:E3').NumberFormat = "@"
$ws.Range('D3').Value = '1.940.91'
$ws.Range('E3').Value = '  -1.30%  '

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5:E5').NumberFormat = "@"
$ws.Range('D5').Value = '243.37'
$ws.Range('E5').Value = '  -1.14%  '

$ws.Range('D6:E6').NumberFormat = "@"
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.13%  '

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.15%  '

$ws.Range('D8:E8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2939'
$ws.Range('E8').Value = '  -1.85%  '

$ws.Range('D9:E9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06906'
$ws.Range('E9').Value = '  +0.14%  '

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.38%  '

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.40%  '

$ws.Range('D12:E12').NumberFormat = "@"
$ws.Range('D12').Value = '1.944.99'
$ws.Range('E12').Value = '  -0.87%  '

$ws.Range('D13:E13').NumberFormat = "@"
$ws.Range('D13').Value = '0.07768'
$ws.Range('E13').Value = '  -0.09%  '

$ws.Range('D14:E14').NumberFormat = "@"
$ws.Range('D14').Value = '5.375'
$ws.Range('E14').Value = '  -2.13%  '

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -1.96%  '

$ws.Range('D16:E16').NumberFormat = "@"
$ws.Range('D16').Value = '275.50'
$ws.Range('E16').Value = '  -4.20%  '

$ws.Range('D17:E17').NumberFormat = "@"
$ws.Range('D17').Value = '30.849.77'
$ws.Range('E17').Value = '  -0.97%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007734'

$ws.Range('D19:E19').NumberFormat = "@"
$ws.Range('D19').Value = '13.11'
$ws.Range('E19').Value = '  -1.68%  '

$ws.Range('D20:E20').NumberFormat = "@"
$ws.Range('D20').Value = '5.610'
$ws.Range('E20').Value = '  +1.07%  '

$ws.Range('D21:E21').NumberFormat = "@"
$ws.Range('D21').Value = '2.194.80'
$ws.Range('E21').Value = '  -0.41%  '

$ws.Range('D22:E22').NumberFormat = "@"
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.07%  '

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.19%  '

$ws.Range('D24:E24').NumberFormat = "@"
$ws.Range('D24').Value = '6.544'
$ws.Range('E24').Value = '  -0.89%  '

$ws.Range('D25:E25').NumberFormat = "@"
$ws.Range('D25').Value = '9.817'
$ws.Range('E25').Value = '  -1.19%  '

$ws.Range('D26:E26').NumberFormat = "@"
$ws.Range('D26').Value = '166.70'
$ws.Range('E26').Value = '  -1.87%  '

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -3.74%  '

$ws.Range('D28:E28').NumberFormat = "@"
$ws.Range('D28').Value = '2.157'
$ws.Range('E28').Value = '  -3.17%  '

$ws.Range('D29:E29').NumberFormat = "@"
$ws.Range('D29').Value = '0.1043'
$ws.Range('E29').Value = '  -1.54%  '

$ws.Range('D30:E30').NumberFormat = "@"
$ws.Range('D30').Value = '1.394'
$ws.Range('E30').Value = '  -3.16%  '

$ws.Range('D31:E31').NumberFormat = "@"
$ws.Range('D31').Value = '1.560'
$ws.Range('E31').Value = '  -1.59%  '

$ws.Range('D32:E32').NumberFormat = "@"
$ws.Range('D32').Value = '4.579'
$ws.Range('E32').Value = '  -1.59%  '

$ws.Range('D33:E33').NumberFormat = "@"
$ws.Range('D33').Value = '4.386'
$ws.Range('E33').Value = '  -2.27%  '

$ws.Range('D34:E34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04889'
$ws.Range('E34').Value = '  -2.22%  '

$ws.Range('D35:E35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7568'
$ws.Range('E35').Value = '  -1.41%  '

$ws.Range('D36:E36').NumberFormat = "@"
$ws.Range('D36').Value = '1.154'
$ws.Range('E36').Value = '  -3.54%  '

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.11%  '

$ws.Range('D38:E38').NumberFormat = "@"
$ws.Range('D38').Value = '2.736'
$ws.Range('E38').Value = '  +0.00%  '

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -2.62%  '

$ws.Range('B40:E40').NumberFormat = "@"
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.655'
$ws.Range('E40').Value = '  -2.02%  '

$ws.Range('B41:E41').NumberFormat = "@"
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '78.69'
$ws.Range('E41').Value = '  +8.40%  '

$ws.Range('D42:E42').NumberFormat = "@"
$ws.Range('D42').Value = '6.479'
$ws.Range('E42').Value = '  +0.44%  '

$ws.Range('D43:E43').NumberFormat = "@"
$ws.Range('D43').Value = '2.093'
$ws.Range('E43').Value = '  -4.96%  '

$ws.Range('D44:E44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9137'
$ws.Range('E44').Value = '  +3.03%  '

$ws.Range('D45:E45').NumberFormat = "@"
$ws.Range('D45').Value = '0.4435'
$ws.Range('E45').Value = '  -2.78%  '

$ws.Range('D46:E46').NumberFormat = "@"
$ws.Range('D46').Value = '107.75'
$ws.Range('E46').Value = '  -1.98%  '

$ws.Range('D47:E47').NumberFormat = "@"
$ws.Range('D47').Value = '0.9987'
$ws.Range('E47').Value = '  -0.37%  '

$ws.Range('D48:E48').NumberFormat = "@"
$ws.Range('D48').Value = '7.684'
$ws.Range('E48').Value = '  -5.68%  '

$ws.Range('D49:E49').NumberFormat = "@"
$ws.Range('D49').Value = '996.83'
$ws.Range('E49').Value = '  +2.81%  '

$ws.Range('D50:E50').NumberFormat = "@"
$ws.Range('D50').Value = '0.1244'
$ws.Range('E50').Value = '  -2.27%  '

$ws.Range('D51:E51').NumberFormat = "@"
$ws.Range('D51').Value = '36.06'
$ws.Range('E51').Value = '  +0.49%  '
